$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Plain text values (column E / F) for the new rows 50-54
# ---------------------------------------------------------------------------
$ws.Range("E50").Value = "사용자 close시 state '대기'로 업데이트"
$ws.Range("F50").Value = "9월12일"

$ws.Range("E51").Value = "state 대기, 완료 -> 대기로 통합"
$ws.Range("F51").Value = "9월12일"

$ws.Range("E52").Value = "분할매수 대기상태 셀값 처리"
$ws.Range("F52").Value = "9월12일"
$ws.Range("G52").Value = "놔둬도됨"

$ws.Range("E53").Value = "물타기 로직"

$ws.Range("E54").Value = "레버리지 제한 오류시 stop"
$ws.Range("F54").Value = "9월12일"

$ws.Range("E56").Value = "테스트"

$ws.Range("E61").Value = "웹에서 매수하면 0차 매수 상태로 변경"

# ---------------------------------------------------------------------------
# 2. The "note" column G/H texts (these get the bold-red font treatment)
# ---------------------------------------------------------------------------
$ws.Range("G53").Value = "최대 물타기 설정값 > ROE 강제 분할매수 -> 소스코드"

$ws.Range("G56").Value = "1.분할매수`r`n○ 1차 분할매수`r`n○ 6차 분할매수"
$ws.Range("G57").Value = "2.익절`r`n○ 1차 익절`r`n○ 6차 익절"
$ws.Range("G58").Value = "3.최대 물타기`r`n○ 1차매수 -> 최대물타기 적용`r`n○ 6차매수 -> 최대물타기 미적용`r`n손절% < 최대물타기%"
$ws.Range("G59").Value = "4.손절`r`n○ 1차매수 -> 손절`r`n○ 2차매수 -> 손절"
$ws.Range("G60").Value = "5.마지막 단계 색상`r`n○ 2단계로 설정 확인"
$ws.Range("G61").Value = "7. 프로그램매매 X, 웹에서 open -> 상태 '대기'`r`n-> [시작] -> 추가매수"

$ws.Range("H56:H61").Value = "설정값:"

# ---------------------------------------------------------------------------
# 3. Formatting: bold + red font on the G53:H61 note block
# ---------------------------------------------------------------------------
$ws.Range("G53:H55").Font.Bold = $true
$ws.Range("G53:H55").Font.Color = 255

$ws.Range("H56:H61").Font.Bold = $true
$ws.Range("H56:H61").Font.Color = 255

$ws.Range("G56:G61").Font.Bold = $true
$ws.Range("G56:G61").Font.Color = 255
$ws.Range("G56:G61").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Row heights for the wrapped, multi-line note rows
# ---------------------------------------------------------------------------
$ws.Rows.Item(56).RowHeight = 49.5
$ws.Rows.Item(57).RowHeight = 49.5
$ws.Rows.Item(58).RowHeight = 66
$ws.Rows.Item(59).RowHeight = 49.5
$ws.Rows.Item(60).RowHeight = 33
$ws.Rows.Item(61).RowHeight = 33

# ---------------------------------------------------------------------------
# 5. Restore selection to match the saved view state
# ---------------------------------------------------------------------------
$ws.Range("G53").Select()

Write-Output "done"
